$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix header typo "Privince" -> "Province" (table column name follows automatically)
$ws.Range("B3").Value = "Province"

# "Chonburi" rows (originally B7:B9) are renamed to "Samut Prakan"
$ws.Range("B7").Value = "Samut Prakan"
$ws.Range("B8").Value = "Samut Prakan"
$ws.Range("B9").Value = "Samut Prakan"

# Update the active selection to match the new state
$ws.Range("B14").Select()
